$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "grouping1"
$ws.Range("E1").Value = "grouping2"

# Match style of existing header cells (A1:C1) so new header cells look consistent
$ws.Range("D1:E1").Font.Color = $ws.Range("A1").Font.Color

# Move active selection to the newly added last header cell, like in the authored file
[void]$ws.Range("E1").Select()
